# This edit rearranges the data rows 2-13 of the "Artfynd" sheet: each
# row's entire record (all columns A:AY) is replaced by the record that,
# in the original workbook, lived on a different row. Concretely the new
# contents come from the old rows via two permutation cycles:
#   2 <- 5 <- 8 <- 10 <- 12 <- 3 <- 6 <- 2
#   4 <- 7 <- 9 <- 11 <- 13 <- 4
# (read "A <- B" as "new row A gets the data that used to be on row B").
#
# We use Copy/PasteSpecial (rather than reading .Value into an array and
# writing it back) because a plain value round-trip lets Excel's "smart"
# text parsing reinterpret date-looking strings (e.g. "2023-03-06") as
# real dates. Copy/PasteSpecial moves the literal cell content instead,
# so text stays text.
#
# Each cycle is applied by stashing the first row of the cycle in a
# scratch row far outside the used range, shifting every other row's data
# one slot "up" the cycle, and finally dropping the stashed data into the
# last slot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = "AY"
$scratchRow = 100

function Copy-RowData($fromRow, $toRow) {
    # Clear the destination first: PasteSpecial leaves a cell untouched
    # when the copied source cell is blank (it doesn't overwrite with a
    # blank), so without this a stale value could survive the paste.
    $ws.Range("A" + $toRow + ":" + $lastCol + $toRow).ClearContents()
    $ws.Range("A" + $fromRow + ":" + $lastCol + $fromRow).Copy()
    $ws.Range("A" + $toRow + ":" + $lastCol + $toRow).PasteSpecial()
}

$cycles = @(
    , @(2, 5, 8, 10, 12, 3, 6)
    , @(4, 7, 9, 11, 13)
)

foreach ($cycle in $cycles) {
    $n = $cycle.Length

    # Stash the original contents of the cycle's first row.
    Copy-RowData $cycle[0] $scratchRow

    # Shift data along the cycle: row[i] <- row[i+1].
    for ($i = 0; $i -lt ($n - 1); $i++) {
        Copy-RowData $cycle[$i + 1] $cycle[$i]
    }

    # Close the loop with the stashed data.
    Copy-RowData $scratchRow $cycle[$n - 1]
}

# Clean up the scratch row so it doesn't linger in the saved workbook.
$ws.Range("A" + $scratchRow + ":" + $lastCol + $scratchRow).ClearContents()
$excel.CutCopyMode = $false
